$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: columns D:T (4-20) get width ~9 characters ---
$ws.Range("D1:T1").EntireColumn.ColumnWidth = 8.14

# --- Copy formatting for the new 2023 column (T) from the 2022 column (S) ---
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)   # xlPasteFormats

# --- New data for year 2023 ---
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 29.5
$ws.Range("T6").Value = 12030.6

# --- Updated 2022 renewable-energy-share figure ---
$ws.Range("S5").Value = 29.9

# --- Footnote label translation fix ---
$ws.Range("C2").Value = "(in percent)"
